$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-11-30 Thursday" "2023-12-01 Friday"

Replace-Text "34×53=1802" "22×95=2090"
Replace-Text "37×65=2405" "71×74=5254"
Replace-Text "70×16=1120" "25×37=925"
Replace-Text "12×70=840" "60×79=4740"
Replace-Text "86×26=2236" "30×36=1080"

Replace-Text "39×82=3198" "72×74=5328"
Replace-Text "15×36=540" "29×85=2465"
Replace-Text "59×82=4838" "16×36=576"
Replace-Text "82×30=2460" "56×81=4536"
Replace-Text "90×26=2340" "80×29=2320"

Replace-Text "71×97=6887" "29×96=2784"
Replace-Text "76×61=4636" "50×86=4300"
Replace-Text "47×36=1692" "83×25=2075"
Replace-Text "22×54=1188" "52×62=3224"
Replace-Text "61×73=4453" "98×85=8330"

Replace-Text "33×46=1518" "32×86=2752"
Replace-Text "90×72=6480" "47×68=3196"
Replace-Text "56×15=840" "43×33=1419"
Replace-Text "13×95=1235" "81×70=5670"
Replace-Text "32×64=2048" "76×31=2356"

Replace-Text "52×45=2340" "41×91=3731"
Replace-Text "33×57=1881" "50×78=3900"
Replace-Text "63×26=1638" "49×74=3626"
Replace-Text "88×91=8008" "24×42=1008"
Replace-Text "74×41=3034" "72×22=1584"
